$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the date line in the first paragraph
$d.Content.Find.Execute("2024-12-06 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-12-07 Saturday", 1) | Out-Null

# Update each division problem cell by exact (row, column) address so that
# duplicate cell texts (e.g. "864÷3=" appearing twice) are each replaced
# with their own, independent target value.
$c = $t.Cell(1, 1)
$c.Range.Find.Execute("770÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "110÷6=", 1) | Out-Null
$c = $t.Cell(1, 2)
$c.Range.Find.Execute("464÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "649÷5=", 1) | Out-Null
$c = $t.Cell(1, 3)
$c.Range.Find.Execute("714÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "540÷5=", 1) | Out-Null
$c = $t.Cell(1, 4)
$c.Range.Find.Execute("407÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "194÷7=", 1) | Out-Null
$c = $t.Cell(1, 5)
$c.Range.Find.Execute("931÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "563÷5=", 1) | Out-Null
$c = $t.Cell(5, 1)
$c.Range.Find.Execute("743÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "223÷7=", 1) | Out-Null
$c = $t.Cell(5, 2)
$c.Range.Find.Execute("621÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "454÷2=", 1) | Out-Null
$c = $t.Cell(5, 3)
$c.Range.Find.Execute("700÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "234÷9=", 1) | Out-Null
$c = $t.Cell(5, 4)
$c.Range.Find.Execute("155÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "407÷9=", 1) | Out-Null
$c = $t.Cell(5, 5)
$c.Range.Find.Execute("692÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "548÷7=", 1) | Out-Null
$c = $t.Cell(9, 1)
$c.Range.Find.Execute("864÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "577÷9=", 1) | Out-Null
$c = $t.Cell(9, 2)
$c.Range.Find.Execute("712÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "643÷2=", 1) | Out-Null
$c = $t.Cell(9, 3)
$c.Range.Find.Execute("291÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "271÷6=", 1) | Out-Null
$c = $t.Cell(9, 4)
$c.Range.Find.Execute("617÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "652÷6=", 1) | Out-Null
$c = $t.Cell(9, 5)
$c.Range.Find.Execute("699÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "245÷5=", 1) | Out-Null
$c = $t.Cell(13, 1)
$c.Range.Find.Execute("864÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "127÷7=", 1) | Out-Null
$c = $t.Cell(13, 2)
$c.Range.Find.Execute("578÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "880÷2=", 1) | Out-Null
$c = $t.Cell(13, 3)
$c.Range.Find.Execute("629÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "410÷8=", 1) | Out-Null
$c = $t.Cell(13, 4)
$c.Range.Find.Execute("105÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "903÷6=", 1) | Out-Null
$c = $t.Cell(13, 5)
$c.Range.Find.Execute("401÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "586÷2=", 1) | Out-Null
$c = $t.Cell(17, 1)
$c.Range.Find.Execute("306÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "803÷2=", 1) | Out-Null
$c = $t.Cell(17, 2)
$c.Range.Find.Execute("685÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "452÷2=", 1) | Out-Null
$c = $t.Cell(17, 3)
$c.Range.Find.Execute("582÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "509÷9=", 1) | Out-Null
$c = $t.Cell(17, 4)
$c.Range.Find.Execute("389÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "864÷8=", 1) | Out-Null
$c = $t.Cell(17, 5)
$c.Range.Find.Execute("856÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "652÷4=", 1) | Out-Null

Write-Output "Done."
